$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row (row 12) describing the "Merchant" class, mirroring the
# pattern of the existing rows (e.g. row 11 / Prisoner).
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Merchant"
$ws.Range("C12").Value = "str"
$ws.Range("D12").Value = "dex"
$ws.Range("P12").Value = "Blacksmith"
$ws.Range("Q12").Value = "Ranger"
$ws.Range("R12").Value = 50
$ws.Range("S12").Value = 30
